# The workbook originally has a single worksheet named "Facilities Maintenance"
# containing a mix of Library and Facilities-Maintenance tender rows.
# This edit splits it into two worksheets:
#   1) "Library"               - renamed from the original sheet, holding the
#                                 Library/NLB/NYP invitation-to-quote rows.
#   2) "Facilities Maintenance" - a newly added sheet (placed after "Library")
#                                 holding the cleaned-up facilities maintenance
#                                 tender rows (POC name/email/tel separated out).

$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing sheet to "Library" and replace its data ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Library"

# Wipe the old combined data (keep header row intact, it's identical) and
# rewrite it with the Library-only rows.
$ws1.Range("A2:K6").ClearContents()

$ws1.Cells.Item(1,1).Value = "Reference Number"
$ws1.Cells.Item(1,2).Value = "Job Description"
$ws1.Cells.Item(1,3).Value = "Agency"
$ws1.Cells.Item(1,4).Value = "Closing Date"
$ws1.Cells.Item(1,5).Value = "Procurement Type"
$ws1.Cells.Item(1,6).Value = "Site Briefing"
$ws1.Cells.Item(1,7).Value = "Compulsary"
$ws1.Cells.Item(1,8).Value = "Briefing Date"
$ws1.Cells.Item(1,9).Value = "POC"
$ws1.Cells.Item(1,10).Value = "POC Email"
$ws1.Cells.Item(1,11).Value = "POC Tel"

$ws1.Cells.Item(2,1).Value = "NLB000ETQ21000035"
$ws1.Cells.Item(2,2).Value = "Invitation to Quote for the Support and Maintenance of Inventorisation System for the National Library Board"
$ws1.Cells.Item(2,3).Value = "National Library Board"
$ws1.Cells.Item(2,4).Value = "14 May 2021`n01:00PM"
$ws1.Cells.Item(2,5).Value = "Invitation to Quote"
$ws1.Cells.Item(2,6).Value = "no"
$ws1.Cells.Item(2,9).Value = "LIN YANFEN"
$ws1.Cells.Item(2,10).Value = "Lin_Yanfen@nlb.gov.sg"

$ws1.Cells.Item(3,1).Value = "NLB000ETQ21000036"
$ws1.Cells.Item(3,2).Value = "INVITATION TO QUOTE FOR THE PROVISION OF MAINTENANCE SERVICES AND ENHANCEMENTS TO THE WEB ARCHIVE SINGAPORE PORTAL FOR NATIONAL LIBRARY BOARD"
$ws1.Cells.Item(3,3).Value = "National Library Board"
$ws1.Cells.Item(3,4).Value = "14 May 2021`n01:00PM"
$ws1.Cells.Item(3,5).Value = "Invitation to Quote"
$ws1.Cells.Item(3,6).Value = "no"
$ws1.Cells.Item(3,9).Value = "CHARLES WIJAYA"
$ws1.Cells.Item(3,10).Value = "Charles_Wijaya@nlb.gov.sg"

$ws1.Cells.Item(4,1).Value = "NYP000ETQ21000121"
$ws1.Cells.Item(4,2).Value = "Jointly operate and train students at the Nanyang Polytechnic Library Training Cafe with NYP School of Business Management"
$ws1.Cells.Item(4,3).Value = "Nanyang Polytechnic"
$ws1.Cells.Item(4,4).Value = "24 May 2021`n01:00PM"
$ws1.Cells.Item(4,5).Value = "Invitation to Quote"
$ws1.Cells.Item(4,6).Value = "Yes"
$ws1.Cells.Item(4,7).Value = "Yes"
$ws1.Cells.Item(4,8).Value = "10 May 2021"
$ws1.Cells.Item(4,9).Value = "LINDA LIM / PATRICK PNG"
$ws1.Cells.Item(4,10).Value = "linda_sf_lim@nyp.edu.sg"

$ws1.Cells.Item(5,1).Value = "NLB000ETQ21000030"
$ws1.Cells.Item(5,2).Value = "INVITATION TO QUOTE FOR THE SUPPLY AND DELIVERY OF 55 INCH TOUCHSCREEN WITH PC AND CASTOR WHEELS FOR THE NATIONAL LIBRARY BOARD"
$ws1.Cells.Item(5,3).Value = "National Library Board"
$ws1.Cells.Item(5,4).Value = "05 May 2021`n01:00PM"
$ws1.Cells.Item(5,5).Value = "Invitation to Quote"
$ws1.Cells.Item(5,6).Value = "no"
$ws1.Cells.Item(5,9).Value = "PRIMARYAZMI SAATazmi_saat@nlb.gov.sglayout_RepaintAllLayouts();6704 1099layout_RepaintAllLayouts();layout_RepaintAllLayouts();NLBlayout_RepaintAllLayouts();"
$ws1.Cells.Item(5,10).Value = "AZMI SAAT"

# --- Step 2: add a new "Facilities Maintenance" sheet after "Library" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Facilities Maintenance"

$ws2.Cells.Item(1,1).Value = "Reference Number"
$ws2.Cells.Item(1,2).Value = "Job Description"
$ws2.Cells.Item(1,3).Value = "Agency"
$ws2.Cells.Item(1,4).Value = "Closing Date"
$ws2.Cells.Item(1,5).Value = "Procurement Type"
$ws2.Cells.Item(1,6).Value = "Site Briefing"
$ws2.Cells.Item(1,7).Value = "Compulsary"
$ws2.Cells.Item(1,8).Value = "Briefing Date"
$ws2.Cells.Item(1,9).Value = "POC"
$ws2.Cells.Item(1,10).Value = "POC Email"
$ws2.Cells.Item(1,11).Value = "POC Tel"

$ws2.Cells.Item(2,1).Value = "JTC000ETT21000013 "
$ws2.Cells.Item(2,2).Value = "APPOINTMENT OF FACILITY MANAGEMENT COMPANY FOR FACILITIES MAINTENANCE OF CENTRAL DEPARTMENT 3 PROPERTIES FOR A PERIOD OF 3 YEARS"
$ws2.Cells.Item(2,3).Value = "Jurong Town Corporation"
$ws2.Cells.Item(2,4).Value = "14 May 2021`n04:00PM"
$ws2.Cells.Item(2,5).Value = "Tender"
$ws2.Cells.Item(2,6).Value = "Yes"
$ws2.Cells.Item(2,7).Value = "Yes"
$ws2.Cells.Item(2,8).Value = "10 Mar 2021"
$ws2.Cells.Item(2,9).Value = "IVY SIM"
$ws2.Cells.Item(2,10).Value = "Ivy_SIM@jtc.gov.sg"

$ws2.Cells.Item(3,1).Value = "NPB000ETT21000039 "
$ws2.Cells.Item(3,2).Value = "TERM CONTRACT FOR THE MAINTENANCE AND UPGRADING WORKS OF PARK FACILITIES FOR PARKS IN PARKS SOUTH WEST BRANCH FOR A PERIOD OF 3 YEARS"
$ws2.Cells.Item(3,3).Value = "National Parks Board"
$ws2.Cells.Item(3,4).Value = "17 May 2021`n04:00PM"
$ws2.Cells.Item(3,5).Value = "Tender"
$ws2.Cells.Item(3,6).Value = "no"
$ws2.Cells.Item(3,9).Value = "ONG SOH HIAN"
$ws2.Cells.Item(3,10).Value = "ONG_SOH_HIAN@NPARKS.GOV.SG"

$ws2.Cells.Item(4,1).Value = "NPB000ETT21000037 "
$ws2.Cells.Item(4,2).Value = "TERM CONTRACT FOR MAINTENANCE AND UPGRADING WORKS FOR PARK FACILITIES IN SOUTH EAST BRANCH FOR A PERIOD OF THREE (3) YEARS"
$ws2.Cells.Item(4,3).Value = "National Parks Board"
$ws2.Cells.Item(4,4).Value = "10 May 2021`n04:00PM"
$ws2.Cells.Item(4,5).Value = "Tender"
$ws2.Cells.Item(4,6).Value = "no"
$ws2.Cells.Item(4,9).Value = "ONG SOH HIAN"
$ws2.Cells.Item(4,10).Value = "ONG_SOH_HIAN@NPARKS.GOV.SG"

$ws2.Cells.Item(5,1).Value = "PAR000ETT21000006 "
$ws2.Cells.Item(5,2).Value = "MAINTENANCE SERVICES FOR BUILDING WORKS AND FACILITIES AT PARLIAMENT HOUSE FOR A PERIOD OF THREE (3) YEARS"
$ws2.Cells.Item(5,3).Value = "Parliament"
$ws2.Cells.Item(5,4).Value = "10 May 2021`n04:00PM"
$ws2.Cells.Item(5,5).Value = "Tender"
$ws2.Cells.Item(5,6).Value = "Yes"
$ws2.Cells.Item(5,7).Value = "Yes"
$ws2.Cells.Item(5,8).Value = "21 Apr 2021"
$ws2.Cells.Item(5,9).Value = "MOHAMMAD AMIN HAMID"
$ws2.Cells.Item(5,10).Value = "Mohammad_Amin_HAMID@parl.gov.sg"

$ws2.Cells.Item(6,1).Value = "URA000ETT21000011 "
$ws2.Cells.Item(6,2).Value = "Replacement Of The Facility Management System & Extra Low Voltage Systems For A Facility At Marina Bay With Option For Maintenance Up To 9 Yrs"
$ws2.Cells.Item(6,3).Value = "Urban Redevelopment Authority"
$ws2.Cells.Item(6,4).Value = "21 May 2021`n04:00PM"
$ws2.Cells.Item(6,5).Value = "Tender"
$ws2.Cells.Item(6,6).Value = "Yes"
$ws2.Cells.Item(6,7).Value = "No"
$ws2.Cells.Item(6,8).Value = "21 Apr 2021"
$ws2.Cells.Item(6,9).Value = "PRIMARYTAN YUAN HONGTAN_Yuan_Hong@ura.gov.sglayout_RepaintAllLayouts();6321 8213layout_RepaintAllLayouts();layout_RepaintAllLayouts();45 Maxwell Road, The URA Centrelayout_RepaintAllLayouts();"
$ws2.Cells.Item(6,10).Value = "TAN YUAN HONG"

# --- Step 3: make "Library" the active/selected sheet (matches activeTab="0") ---
$ws1.Activate()
